$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.991.18'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.753.85'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '336.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3837'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3420'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.94'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.93%  '
$ws.Range("E10").Value = '  -1.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07237'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.78%  '
$ws.Range("E12").Value = '  +0.85%  '
$ws.Range("E13").Value = '  -0.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.172'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.152'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.750.04'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.70%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001062'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06607'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("E19").Value = '  -3.31%  '
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.200'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.88%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.989.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.62%  '
$ws.Range("E24").Value = '  -3.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.375'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.17'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.88'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.302'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.950.29'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.263'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -11.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '130.98'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.029'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.852'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08818'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.23'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.556'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6576'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.53%  '
$ws.Range("E38").Value = '  -4.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.154'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06159'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2107'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.215'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.979'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9989'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.81'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.840'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6060'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.20'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.15%  '
$ws.Range("E49").Value = '  -2.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.167'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.114'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.96%  '
